# Generate Report for Handback
#
# The localization handback process regenerated the XLIFF handoff/handback
# report, producing fresh "Latest Handoff Datetime" / "Latest Handback
# DateTime" timestamps for the file that was just handed back
# (6055ed7f-39dc-4826-8695-268d6bc539a8) on both the "zh-cn" and "de-de"
# per-language status sheets.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 holds 6055ed7f-39dc-4826-8695-268d6bc539a8
#   H = Latest Handoff Datetime, K = Latest Handback DateTime
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("H2").Value = "2016-08-18 10:51:27"
$wsZh.Range("K2").Value = "2016-08-18 10:51:42"

# de-de sheet: row 2 holds 6055ed7f-39dc-4826-8695-268d6bc539a8
#   H = Latest Handoff Datetime, K = Latest Handback DateTime
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("H2").Value = "2016-08-18 10:51:32"
$wsDe.Range("K2").Value = "2016-08-18 10:51:50"
